$wb = $excel.ActiveWorkbook

# ============================================================
# Sheet "Overview": insert a new row 2 for the new handed-off
# file (15cf328a-...), pushing the existing 42ce9edd-... row
# down to row 3.
# ============================================================
$ws = $wb.Worksheets.Item("Overview")
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "15cf328a-2d92-4bdc-b1e2-3d721be47214.md"
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-44-13 22:44:31"

$ws.Range("A3").Value = "42ce9edd-5f08-418e-a340-64200725e749.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-43-13 22:43:31"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/34e3ec46bc5e6ad2fd438750a032bd2a0446b950/e2e/15cf328a-2d92-4bdc-b1e2-3d721be47214.md", "", "", "15cf328a-2d92-4bdc-b1e2-3d721be47214.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/34e3ec46bc5e6ad2fd438750a032bd2a0446b950/e2e/42ce9edd-5f08-418e-a340-64200725e749.md", "", "", "42ce9edd-5f08-418e-a340-64200725e749.md") | Out-Null

# ============================================================
# Sheet "zh-cn": insert a new row 2 for the new handed-off file,
# pushing the existing row down to row 3.
# ============================================================
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "15cf328a-2d92-4bdc-b1e2-3d721be47214.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "15cf328a-2d92-4bdc-b1e2-3d721be47214.1ff711125536c715d1211fb1a613192d8afb6e82.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-13 22:44:27"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

$ws.Range("A3").Value = "42ce9edd-5f08-418e-a340-64200725e749.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "42ce9edd-5f08-418e-a340-64200725e749.cbf6ae567a532461c0a8c6f0730fe79c0bf1b33c.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-13 22:40:56"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/34e3ec46bc5e6ad2fd438750a032bd2a0446b950/e2e/15cf328a-2d92-4bdc-b1e2-3d721be47214.md", "", "", "15cf328a-2d92-4bdc-b1e2-3d721be47214.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/34e3ec46bc5e6ad2fd438750a032bd2a0446b950/e2e/15cf328a-2d92-4bdc-b1e2-3d721be47214.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9ba253d6763ff74ea83c3b8dfedeaa494f35726b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/15cf328a-2d92-4bdc-b1e2-3d721be47214.1ff711125536c715d1211fb1a613192d8afb6e82.zh-cn.xlf", "", "", "15cf328a-2d92-4bdc-b1e2-3d721be47214.1ff711125536c715d1211fb1a613192d8afb6e82.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/34e3ec46bc5e6ad2fd438750a032bd2a0446b950/e2e/42ce9edd-5f08-418e-a340-64200725e749.md", "", "", "42ce9edd-5f08-418e-a340-64200725e749.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/34e3ec46bc5e6ad2fd438750a032bd2a0446b950/e2e/42ce9edd-5f08-418e-a340-64200725e749.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9ba253d6763ff74ea83c3b8dfedeaa494f35726b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/42ce9edd-5f08-418e-a340-64200725e749.cbf6ae567a532461c0a8c6f0730fe79c0bf1b33c.zh-cn.xlf", "", "", "42ce9edd-5f08-418e-a340-64200725e749.cbf6ae567a532461c0a8c6f0730fe79c0bf1b33c.zh-cn.xlf") | Out-Null

$ws.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ============================================================
# Sheet "de-de": insert a new row 2 for the new handed-off file,
# pushing the existing row down to row 3.
# ============================================================
$ws = $wb.Worksheets.Item("de-de")
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "15cf328a-2d92-4bdc-b1e2-3d721be47214.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "15cf328a-2d92-4bdc-b1e2-3d721be47214.1ff711125536c715d1211fb1a613192d8afb6e82.de-de.xlf"
$ws.Range("E2").Value = "2016-03-13 22:44:31"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

$ws.Range("A3").Value = "42ce9edd-5f08-418e-a340-64200725e749.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "42ce9edd-5f08-418e-a340-64200725e749.cbf6ae567a532461c0a8c6f0730fe79c0bf1b33c.de-de.xlf"
$ws.Range("E3").Value = "2016-03-13 22:43:31"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/34e3ec46bc5e6ad2fd438750a032bd2a0446b950/e2e/15cf328a-2d92-4bdc-b1e2-3d721be47214.md", "", "", "15cf328a-2d92-4bdc-b1e2-3d721be47214.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/34e3ec46bc5e6ad2fd438750a032bd2a0446b950/e2e/15cf328a-2d92-4bdc-b1e2-3d721be47214.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/41512d59c850a277610f71484717facd368d5608/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/15cf328a-2d92-4bdc-b1e2-3d721be47214.1ff711125536c715d1211fb1a613192d8afb6e82.de-de.xlf", "", "", "15cf328a-2d92-4bdc-b1e2-3d721be47214.1ff711125536c715d1211fb1a613192d8afb6e82.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/34e3ec46bc5e6ad2fd438750a032bd2a0446b950/e2e/42ce9edd-5f08-418e-a340-64200725e749.md", "", "", "42ce9edd-5f08-418e-a340-64200725e749.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/34e3ec46bc5e6ad2fd438750a032bd2a0446b950/e2e/42ce9edd-5f08-418e-a340-64200725e749.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/41512d59c850a277610f71484717facd368d5608/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/42ce9edd-5f08-418e-a340-64200725e749.cbf6ae567a532461c0a8c6f0730fe79c0bf1b33c.de-de.xlf", "", "", "42ce9edd-5f08-418e-a340-64200725e749.cbf6ae567a532461c0a8c6f0730fe79c0bf1b33c.de-de.xlf") | Out-Null

$ws.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

Write-Host "Report for handoff generated."
